$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.85292862532423
$ws.Range("C2").Value = 7.994163364222261
$ws.Range("D2").Value = 7.484255606184355
$ws.Range("F2").Value = 40.49427171907251
$ws.Range("G2").Value = 48.01132610198719
$ws.Range("H2").Value = 18.91783243392204
$ws.Range("J2").Value = 10.31135072363204
$ws.Range("L2").Value = 11.97902579733747
$ws.Range("N2").Value = 19.28708406754004
$ws.Range("B3").Value = 21.46850249411604
$ws.Range("C3").Value = 7.628238563502476
$ws.Range("D3").Value = 7.486012396679017
$ws.Range("F3").Value = 40.49690934936429
$ws.Range("G3").Value = 47.90370212990973
$ws.Range("H3").Value = 18.95464882032681
$ws.Range("J3").Value = 10.33535329139253
$ws.Range("L3").Value = 11.97215548735251
$ws.Range("N3").Value = 19.35930454826159
$ws.Range("B4").Value = 21.23393765387704
$ws.Range("C4").Value = 7.39284617914939
$ws.Range("D4").Value = 7.487535840291744
$ws.Range("F4").Value = 40.50916192268015
$ws.Range("G4").Value = 47.85305341749289
$ws.Range("H4").Value = 18.98122879686877
$ws.Range("J4").Value = 10.35092376328262
$ws.Range("L4").Value = 11.96969354966122
$ws.Range("N4").Value = 19.4056207432084
$ws.Range("B5").Value = 21.13885115210162
$ws.Range("C5").Value = 7.294298418682835
$ws.Range("D5").Value = 7.488268909409356
$ws.Range("F5").Value = 40.51682473634794
$ws.Range("G5").Value = 47.83630136146149
$ws.Range("H5").Value = 18.99305732058004
$ws.Range("J5").Value = 10.35747881540272
$ws.Range("L5").Value = 11.96913351139343
$ws.Range("N5").Value = 19.42499273515961
$ws.Range("B6").Value = 21.12309615178357
$ws.Range("C6").Value = 7.2777784790724
$ws.Range("D6").Value = 7.488397426790098
$ws.Range("F6").Value = 40.51825824035782
$ws.Range("G6").Value = 47.83375461501348
$ws.Range("H6").Value = 18.99508157616278
$ws.Range("J6").Value = 10.3585799738306
$ws.Range("L6").Value = 11.96906732840785
$ws.Range("N6").Value = 19.42823955518935
$ws.Range("B7").Value = 21.23265308524129
$ws.Range("C7").Value = 7.391527650115276
$ws.Range("D7").Value = 7.487545271653695
$ws.Range("F7").Value = 40.50925446279751
$ws.Range("G7").Value = 47.85281174664172
$ws.Range("H7").Value = 18.98138428726432
$ws.Range("J7").Value = 10.35101131617256
$ws.Range("L7").Value = 11.96968420030205
$ws.Range("N7").Value = 19.4058799829287
$ws.Range("B8").Value = 21.72015008911644
$ws.Range("C8").Value = 7.87025471558983
$ws.Range("D8").Value = 7.484769284742088
$ws.Range("F8").Value = 40.49297238122491
$ws.Range("G8").Value = 47.97101987654467
$ws.Range("H8").Value = 18.92970037475853
$ws.Range("J8").Value = 10.3194542897528
$ws.Range("L8").Value = 11.97629334030212
$ws.Range("N8").Value = 19.31157735455237
$ws.Range("B9").Value = 22.68219294923007
$ws.Range("C9").Value = 8.721494749432217
$ws.Range("D9").Value = 7.48283675129687
$ws.Range("F9").Value = 40.5455193431086
$ws.Range("G9").Value = 48.32475674044171
$ws.Range("H9").Value = 18.85998814422722
$ws.Range("J9").Value = 10.26415414474674
$ws.Range("L9").Value = 12.00311634249159
$ws.Range("N9").Value = 19.14221978286573
$ws.Range("B10").Value = 23.38532976529134
$ws.Range("C10").Value = 9.290800686893972
$ws.Range("D10").Value = 7.483533852108522
$ws.Range("F10").Value = 40.63566842351061
$ws.Range("G10").Value = 48.65797207564802
$ws.Range("H10").Value = 18.82818971667728
$ws.Range("J10").Value = 10.22750373340361
$ws.Range("L10").Value = 12.03116529058619
$ws.Range("N10").Value = 19.02716844041772
$ws.Range("B11").Value = 23.70291809708609
$ws.Range("C11").Value = 9.537193298968379
$ws.Range("D11").Value = 7.484305541907534
$ws.Range("F11").Value = 40.68784376071194
$ws.Range("G11").Value = 48.82516955991784
$ws.Range("H11").Value = 18.81796386605319
$ws.Range("J11").Value = 10.21168710210232
$ws.Range("L11").Value = 12.04570980474427
$ws.Range("N11").Value = 18.97683960725836
$ws.Range("B12").Value = 23.82273615772866
$ws.Range("C12").Value = 9.628659023089829
$ws.Range("D12").Value = 7.484662673358895
$ws.Range("F12").Value = 40.70920127853398
$ws.Range("G12").Value = 48.89069491694321
$ws.Range("H12").Value = 18.81470269488785
$ws.Range("J12").Value = 10.20582027269386
$ws.Range("L12").Value = 12.05147144639918
$ws.Range("N12").Value = 18.95806836605235
$ws.Range("B13").Value = 23.79695285421898
$ws.Range("C13").Value = 9.609042365092154
$ws.Range("D13").Value = 7.484582879235933
$ws.Range("F13").Value = 40.70453051376972
$ws.Range("G13").Value = 48.87648510285773
$ws.Range("H13").Value = 18.81537784467818
$ws.Range("J13").Value = 10.2070783549485
$ws.Range("L13").Value = 12.05021932463633
$ws.Range("N13").Value = 18.96209834157882
$ws.Range("B14").Value = 23.71278519617141
$ws.Range("C14").Value = 9.544755211482265
$ws.Range("D14").Value = 7.484333625100153
$ws.Range("F14").Value = 40.68956883318567
$ws.Range("G14").Value = 48.83051626459864
$ws.Range("H14").Value = 18.81768331172225
$ws.Range("J14").Value = 10.21120198045292
$ws.Range("L14").Value = 12.04617874519648
$ws.Range("N14").Value = 18.97528953947013
$ws.Range("B15").Value = 23.66116858944988
$ws.Range("C15").Value = 9.505137404974279
$ws.Range("D15").Value = 7.48418938912003
$ws.Range("F15").Value = 40.68061251476757
$ws.Range("G15").Value = 48.80264590826841
$ws.Range("H15").Value = 18.81917510529473
$ws.Range("J15").Value = 10.21374377022735
$ws.Range("L15").Value = 12.04373676140124
$ws.Range("N15").Value = 18.98340688195433
$ws.Range("B16").Value = 23.36451859885146
$ws.Range("C16").Value = 9.274442560728549
$ws.Range("D16").Value = 7.48349252860326
$ws.Range("F16").Value = 40.63248287806812
$ws.Range("G16").Value = 48.64735673310399
$ws.Range("H16").Value = 18.82894344332652
$ws.Range("J16").Value = 10.22855456759479
$ws.Range("L16").Value = 12.03025046142129
$ws.Range("N16").Value = 19.03049782588295
$ws.Range("B17").Value = 23.18186506135493
$ws.Range("C17").Value = 9.129673276727152
$ws.Range("D17").Value = 7.483181131436377
$ws.Range("F17").Value = 40.60581326156039
$ws.Range("G17").Value = 48.55606853141916
$ws.Range("H17").Value = 18.83602297814923
$ws.Range("J17").Value = 10.23785936527525
$ws.Range("L17").Value = 12.02243233923945
$ws.Range("N17").Value = 19.05989987607732
$ws.Range("B18").Value = 23.0766004968002
$ws.Range("C18").Value = 9.045222296739265
$ws.Range("D18").Value = 7.48304483868966
$ws.Range("F18").Value = 40.59152507494503
$ws.Range("G18").Value = 48.50503417680225
$ws.Range("H18").Value = 18.84049389901901
$ws.Range("J18").Value = 10.24329182148533
$ws.Range("L18").Value = 12.01810374117927
$ws.Range("N18").Value = 19.07700031396131
$ws.Range("B19").Value = 23.04092787469073
$ws.Range("C19").Value = 9.016426311613499
$ws.Range("D19").Value = 7.483006059596701
$ws.Range("F19").Value = 40.58686807797596
$ws.Range("G19").Value = 48.48800864622441
$ws.Range("H19").Value = 18.84207614723022
$ws.Range("J19").Value = 10.24514501270711
$ws.Range("L19").Value = 12.0166671167511
$ws.Range("N19").Value = 19.08282276623574
$ws.Range("B20").Value = 23.20133115426925
$ws.Range("C20").Value = 9.145206921731431
$ws.Range("D20").Value = 7.483209852320638
$ws.Range("F20").Value = 40.60854350090204
$ws.Range("G20").Value = 48.56563418354919
$ws.Range("H20").Value = 18.83522804710259
$ws.Range("J20").Value = 10.2368605171952
$ws.Range("L20").Value = 12.02324720488529
$ws.Range("N20").Value = 19.05675041266957
$ws.Range("B21").Value = 23.73752029560973
$ws.Range("C21").Value = 9.563687979049027
$ws.Range("D21").Value = 7.484405079153234
$ws.Range("F21").Value = 40.69392007765246
$ws.Range("G21").Value = 48.84395870747048
$ws.Range("H21").Value = 18.81698954333692
$ws.Range("J21").Value = 10.20998744855101
$ws.Range("L21").Value = 12.04735869182134
$ws.Range("N21").Value = 18.97140718348935
$ws.Range("B22").Value = 24.08530159289798
$ws.Range("C22").Value = 9.826468126993921
$ws.Range("D22").Value = 7.485564351788416
$ws.Range("F22").Value = 40.75903973375319
$ws.Range("G22").Value = 49.03873008812476
$ws.Range("H22").Value = 18.80863231265453
$ws.Range("J22").Value = 10.1931386482191
$ws.Range("L22").Value = 12.06459575898214
$ws.Range("N22").Value = 18.91730369244559
$ws.Range("B23").Value = 23.89996442425177
$ws.Range("C23").Value = 9.687206242254289
$ws.Range("D23").Value = 7.484911181265158
$ws.Range("F23").Value = 40.72343367081949
$ws.Range("G23").Value = 48.93361175044098
$ws.Range("H23").Value = 18.81276630644343
$ws.Range("J23").Value = 10.20206596868241
$ws.Range("L23").Value = 12.05526164791702
$ws.Range("N23").Value = 18.94602718083191
$ws.Range("B24").Value = 23.19253131759076
$ws.Range("C24").Value = 9.138187962455808
$ws.Range("D24").Value = 7.483196734475957
$ws.Range("F24").Value = 40.60730590566096
$ws.Range("G24").Value = 48.56130503791289
$ws.Range("H24").Value = 18.83558618706978
$ws.Range("J24").Value = 10.23731183782819
$ws.Range("L24").Value = 12.02287828641187
$ws.Range("N24").Value = 19.05817367189938
$ws.Range("B25").Value = 22.4220600708292
$ws.Range("C25").Value = 8.500898445135181
$ws.Range("D25").Value = 7.482986062346845
$ws.Range("F25").Value = 40.52224853926894
$ws.Range("G25").Value = 48.2160932555059
$ws.Range("H25").Value = 18.87544548641587
$ws.Range("J25").Value = 10.27841313380445
$ws.Range("L25").Value = 11.99438638876183
$ws.Range("N25").Value = 19.18638068807378
